# Update TPM-derived NATMI ligand/receptor metrics on Sheet1 (Sema3b-Nrp1)
# with newly recomputed values (per commit "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 7).Value2 = 3.83668
$ws.Cells.Item(2, 8).Value2 = 11.51004
$ws.Cells.Item(2, 9).Value2 = 0.1710971228178894
$ws.Cells.Item(2, 10).Value2 = 0.1710971228178894
$ws.Cells.Item(2, 13).Value2 = 133.7780026666667
$ws.Cells.Item(2, 14).Value2 = 401.334008
$ws.Cells.Item(2, 15).Value2 = 0.50863533211804
$ws.Cells.Item(2, 16).Value2 = 0.5086353321180399
$ws.Cells.Item(2, 17).Value2 = 513.2633872711467
$ws.Cells.Item(2, 18).Value2 = 4619.37048544032
$ws.Cells.Item(2, 19).Value2 = 0.08702604188891824
$ws.Cells.Item(2, 20).Value2 = 0.08702604188891822
$ws.Cells.Item(3, 7).Value2 = 3.83668
$ws.Cells.Item(3, 8).Value2 = 11.51004
$ws.Cells.Item(3, 9).Value2 = 0.1710971228178894
$ws.Cells.Item(3, 10).Value2 = 0.1710971228178894
$ws.Cells.Item(3, 15).Value2 = 0.1993888292903622
$ws.Cells.Item(3, 16).Value2 = 0.1993888292903622
$ws.Cells.Item(3, 17).Value2 = 201.20306129628
$ws.Cells.Item(3, 18).Value2 = 1810.82755166652
$ws.Cells.Item(3, 19).Value2 = 0.03411485501360828
$ws.Cells.Item(3, 20).Value2 = 0.03411485501360827
$ws.Cells.Item(4, 7).Value2 = 3.83668
$ws.Cells.Item(4, 8).Value2 = 11.51004
$ws.Cells.Item(4, 9).Value2 = 0.1710971228178894
$ws.Cells.Item(4, 10).Value2 = 0.1710971228178894
$ws.Cells.Item(4, 13).Value2 = 21.197691
$ws.Cells.Item(4, 14).Value2 = 63.593073
$ws.Cells.Item(4, 15).Value2 = 0.08059542216956049
$ws.Cells.Item(4, 16).Value2 = 0.08059542216956046
$ws.Cells.Item(4, 17).Value2 = 81.32875710588
$ws.Cells.Item(4, 18).Value2 = 731.9588139529201
$ws.Cells.Item(4, 19).Value2 = 0.01378964484550493
$ws.Cells.Item(4, 20).Value2 = 0.01378964484550493
$ws.Cells.Item(5, 7).Value2 = 3.83668
$ws.Cells.Item(5, 8).Value2 = 11.51004
$ws.Cells.Item(5, 9).Value2 = 0.1710971228178894
$ws.Cells.Item(5, 10).Value2 = 0.1710971228178894
$ws.Cells.Item(5, 13).Value2 = 55.59592133333333
$ws.Cells.Item(5, 14).Value2 = 166.787764
$ws.Cells.Item(5, 15).Value2 = 0.2113804164220374
$ws.Cells.Item(5, 16).Value2 = 0.2113804164220373
$ws.Cells.Item(5, 17).Value2 = 213.3037594611733
$ws.Cells.Item(5, 18).Value2 = 1919.73383515056
$ws.Cells.Item(5, 19).Value2 = 0.03616658106985792
$ws.Cells.Item(5, 20).Value2 = 0.03616658106985791
$ws.Cells.Item(6, 9).Value2 = 0.5304251323586242
$ws.Cells.Item(6, 10).Value2 = 0.5304251323586241
$ws.Cells.Item(6, 13).Value2 = 133.7780026666667
$ws.Cells.Item(6, 14).Value2 = 401.334008
$ws.Cells.Item(6, 15).Value2 = 0.50863533211804
$ws.Cells.Item(6, 16).Value2 = 0.5086353321180399
$ws.Cells.Item(6, 17).Value2 = 1591.18865147666
$ws.Cells.Item(6, 18).Value2 = 14320.69786328994
$ws.Cells.Item(6, 19).Value2 = 0.2697929633609842
$ws.Cells.Item(6, 20).Value2 = 0.269792963360984
$ws.Cells.Item(7, 9).Value2 = 0.5304251323586242
$ws.Cells.Item(7, 10).Value2 = 0.5304251323586241
$ws.Cells.Item(7, 15).Value2 = 0.1993888292903622
$ws.Cells.Item(7, 16).Value2 = 0.1993888292903622
$ws.Cells.Item(7, 17).Value2 = 623.757773721494
$ws.Cells.Item(7, 18).Value2 = 5613.819963493446
$ws.Cells.Item(7, 19).Value2 = 0.1057608461671715
$ws.Cells.Item(7, 20).Value2 = 0.1057608461671715
$ws.Cells.Item(8, 9).Value2 = 0.5304251323586242
$ws.Cells.Item(8, 10).Value2 = 0.5304251323586241
$ws.Cells.Item(8, 13).Value2 = 21.197691
$ws.Cells.Item(8, 14).Value2 = 63.593073
$ws.Cells.Item(8, 15).Value2 = 0.08059542216956049
$ws.Cells.Item(8, 16).Value2 = 0.08059542216956046
$ws.Cells.Item(8, 17).Value2 = 252.130579649574
$ws.Cells.Item(8, 18).Value2 = 2269.175216846166
$ws.Cells.Item(8, 19).Value2 = 0.04274983747178832
$ws.Cells.Item(8, 20).Value2 = 0.04274983747178829
$ws.Cells.Item(9, 9).Value2 = 0.5304251323586242
$ws.Cells.Item(9, 10).Value2 = 0.5304251323586241
$ws.Cells.Item(9, 13).Value2 = 55.59592133333333
$ws.Cells.Item(9, 14).Value2 = 166.787764
$ws.Cells.Item(9, 15).Value2 = 0.2113804164220374
$ws.Cells.Item(9, 16).Value2 = 0.2113804164220373
$ws.Cells.Item(9, 17).Value2 = 661.2716390632097
$ws.Cells.Item(9, 18).Value2 = 5951.444751568887
$ws.Cells.Item(9, 19).Value2 = 0.1121214853586803
$ws.Cells.Item(9, 20).Value2 = 0.1121214853586802
$ws.Cells.Item(10, 7).Value2 = 5.772967666666666
$ws.Cells.Item(10, 8).Value2 = 17.318903
$ws.Cells.Item(10, 9).Value2 = 0.2574460621911056
$ws.Cells.Item(10, 10).Value2 = 0.2574460621911055
$ws.Cells.Item(10, 13).Value2 = 133.7780026666667
$ws.Cells.Item(10, 14).Value2 = 401.334008
$ws.Cells.Item(10, 15).Value2 = 0.50863533211804
$ws.Cells.Item(10, 16).Value2 = 0.5086353321180399
$ws.Cells.Item(10, 17).Value2 = 772.2960839059137
$ws.Cells.Item(10, 18).Value2 = 6950.664755153224
$ws.Cells.Item(10, 19).Value2 = 0.1309461633450545
$ws.Cells.Item(10, 20).Value2 = 0.1309461633450545
$ws.Cells.Item(11, 7).Value2 = 5.772967666666666
$ws.Cells.Item(11, 8).Value2 = 17.318903
$ws.Cells.Item(11, 9).Value2 = 0.2574460621911056
$ws.Cells.Item(11, 10).Value2 = 0.2574460621911055
$ws.Cells.Item(11, 15).Value2 = 0.1993888292903622
$ws.Cells.Item(11, 16).Value2 = 0.1993888292903622
$ws.Cells.Item(11, 17).Value2 = 302.745802959271
$ws.Cells.Item(11, 18).Value2 = 2724.712226633439
$ws.Cells.Item(11, 19).Value2 = 0.05133186894569832
$ws.Cells.Item(11, 20).Value2 = 0.0513318689456983
$ws.Cells.Item(12, 7).Value2 = 5.772967666666666
$ws.Cells.Item(12, 8).Value2 = 17.318903
$ws.Cells.Item(12, 9).Value2 = 0.2574460621911056
$ws.Cells.Item(12, 10).Value2 = 0.2574460621911055
$ws.Cells.Item(12, 13).Value2 = 21.197691
$ws.Cells.Item(12, 14).Value2 = 63.593073
$ws.Cells.Item(12, 15).Value2 = 0.08059542216956049
$ws.Cells.Item(12, 16).Value2 = 0.08059542216956046
$ws.Cells.Item(12, 17).Value2 = 122.373584750991
$ws.Cells.Item(12, 18).Value2 = 1101.362262758919
$ws.Cells.Item(12, 19).Value2 = 0.02074897406818308
$ws.Cells.Item(12, 20).Value2 = 0.02074897406818307
$ws.Cells.Item(13, 7).Value2 = 5.772967666666666
$ws.Cells.Item(13, 8).Value2 = 17.318903
$ws.Cells.Item(13, 9).Value2 = 0.2574460621911056
$ws.Cells.Item(13, 10).Value2 = 0.2574460621911055
$ws.Cells.Item(13, 13).Value2 = 55.59592133333333
$ws.Cells.Item(13, 14).Value2 = 166.787764
$ws.Cells.Item(13, 15).Value2 = 0.2113804164220374
$ws.Cells.Item(13, 16).Value2 = 0.2113804164220373
$ws.Cells.Item(13, 17).Value2 = 320.9534562558769
$ws.Cells.Item(13, 18).Value2 = 2888.581106302891
$ws.Cells.Item(13, 19).Value2 = 0.05441905583216962
$ws.Cells.Item(13, 20).Value2 = 0.0544190558321696
$ws.Cells.Item(14, 7).Value2 = 0.9200940000000001
$ws.Cells.Item(14, 8).Value2 = 2.760282
$ws.Cells.Item(14, 9).Value2 = 0.04103168263238089
$ws.Cells.Item(14, 10).Value2 = 0.04103168263238088
$ws.Cells.Item(14, 13).Value2 = 133.7780026666667
$ws.Cells.Item(14, 14).Value2 = 401.334008
$ws.Cells.Item(14, 15).Value2 = 0.50863533211804
$ws.Cells.Item(14, 16).Value2 = 0.5086353321180399
$ws.Cells.Item(14, 17).Value2 = 123.088337585584
$ws.Cells.Item(14, 18).Value2 = 1107.795038270256
$ws.Cells.Item(14, 19).Value2 = 0.02087016352308307
$ws.Cells.Item(14, 20).Value2 = 0.02087016352308306
$ws.Cells.Item(15, 7).Value2 = 0.9200940000000001
$ws.Cells.Item(15, 8).Value2 = 2.760282
$ws.Cells.Item(15, 9).Value2 = 0.04103168263238089
$ws.Cells.Item(15, 10).Value2 = 0.04103168263238088
$ws.Cells.Item(15, 15).Value2 = 0.1993888292903622
$ws.Cells.Item(15, 16).Value2 = 0.1993888292903622
$ws.Cells.Item(15, 17).Value2 = 48.251542865274
$ws.Cells.Item(15, 18).Value2 = 434.2638857874661
$ws.Cells.Item(15, 19).Value2 = 0.008181259163884114
$ws.Cells.Item(15, 20).Value2 = 0.008181259163884111
$ws.Cells.Item(16, 7).Value2 = 0.9200940000000001
$ws.Cells.Item(16, 8).Value2 = 2.760282
$ws.Cells.Item(16, 9).Value2 = 0.04103168263238089
$ws.Cells.Item(16, 10).Value2 = 0.04103168263238088
$ws.Cells.Item(16, 13).Value2 = 21.197691
$ws.Cells.Item(16, 14).Value2 = 63.593073
$ws.Cells.Item(16, 15).Value2 = 0.08059542216956049
$ws.Cells.Item(16, 16).Value2 = 0.08059542216956046
$ws.Cells.Item(16, 17).Value2 = 19.503868302954
$ws.Cells.Item(16, 18).Value2 = 175.534814726586
$ws.Cells.Item(16, 19).Value2 = 0.003306965784084161
$ws.Cells.Item(16, 20).Value2 = 0.003306965784084159
$ws.Cells.Item(17, 7).Value2 = 0.9200940000000001
$ws.Cells.Item(17, 8).Value2 = 2.760282
$ws.Cells.Item(17, 9).Value2 = 0.04103168263238089
$ws.Cells.Item(17, 10).Value2 = 0.04103168263238088
$ws.Cells.Item(17, 13).Value2 = 55.59592133333333
$ws.Cells.Item(17, 14).Value2 = 166.787764
$ws.Cells.Item(17, 15).Value2 = 0.2113804164220374
$ws.Cells.Item(17, 16).Value2 = 0.2113804164220373
$ws.Cells.Item(17, 17).Value2 = 51.153473643272
$ws.Cells.Item(17, 18).Value2 = 460.381262789448
$ws.Cells.Item(17, 19).Value2 = 0.00867329416132955
$ws.Cells.Item(17, 20).Value2 = 0.008673294161329547
